$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date in column C for rows 2 through 16
# from 45233 (2023-11-03) to 45243 (2023-11-13)
for ($row = 2; $row -le 16; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45233) {
        $cell.Value = 45243
    }
}
